$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.6742976666666666
$ws.Range("H2").Value = 2.022893
$ws.Range("I2").Value = 0.3960292783506769
$ws.Range("J2").Value = 0.3960292783506769
$ws.Range("M2").Value = 8.432170666666666
$ws.Range("N2").Value = 25.296512
$ws.Range("O2").Value = 0.0153412147997323
$ws.Range("P2").Value = 0.01534121479973231
$ws.Range("Q2").Value = 5.685793005468443
$ws.Range("R2").Value = 51.172137049216
$ws.Range("S2").Value = 0.006075570226160709
$ws.Range("T2").Value = 0.00607557022616071
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.6742976666666666
$ws.Range("H3").Value = 2.022893
$ws.Range("I3").Value = 0.3960292783506769
$ws.Range("J3").Value = 0.3960292783506769
$ws.Range("M3").Value = 211.5004576666667
$ws.Range("N3").Value = 634.5013730000001
$ws.Range("O3").Value = 0.3847969970689267
$ws.Range("P3").Value = 0.3847969970689267
$ws.Range("Q3").Value = 142.6142651035654
$ws.Range("R3").Value = 1283.528385932089
$ws.Range("S3").Value = 0.1523908770607146
$ws.Range("T3").Value = 0.1523908770607146
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.6742976666666666
$ws.Range("H4").Value = 2.022893
$ws.Range("I4").Value = 0.3960292783506769
$ws.Range("J4").Value = 0.3960292783506769
$ws.Range("M4").Value = 149.6042426666667
$ws.Range("N4").Value = 448.812728
$ws.Range("O4").Value = 0.2721850532240109
$ws.Range("P4").Value = 0.2721850532240109
$ws.Range("Q4").Value = 100.8777917535671
$ws.Range("R4").Value = 907.9001257821039
$ws.Range("S4").Value = 0.1077932502061456
$ws.Range("T4").Value = 0.1077932502061456
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.6742976666666666
$ws.Range("H5").Value = 2.022893
$ws.Range("I5").Value = 0.3960292783506769
$ws.Range("J5").Value = 0.3960292783506769
$ws.Range("M5").Value = 180.1047823333333
$ws.Range("N5").Value = 540.314347
$ws.Range("O5").Value = 0.3276767349073302
$ws.Range("P5").Value = 0.3276767349073302
$ws.Range("Q5").Value = 121.4442344828745
$ws.Range("R5").Value = 1092.998110345871
$ws.Range("S5").Value = 0.129769580857656
$ws.Range("T5").Value = 0.129769580857656
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.8625470000000001
$ws.Range("H6").Value = 2.587641
$ws.Range("I6").Value = 0.5065920925430184
$ws.Range("J6").Value = 0.5065920925430184
$ws.Range("M6").Value = 8.432170666666666
$ws.Range("N6").Value = 25.296512
$ws.Range("O6").Value = 0.0153412147997323
$ws.Range("P6").Value = 0.01534121479973231
$ws.Range("Q6").Value = 7.273143512021333
$ws.Range("R6").Value = 65.458291608192
$ws.Range("S6").Value = 0.007771738107548311
$ws.Range("T6").Value = 0.007771738107548312
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.8625470000000001
$ws.Range("H7").Value = 2.587641
$ws.Range("I7").Value = 0.5065920925430184
$ws.Range("J7").Value = 0.5065920925430184
$ws.Range("M7").Value = 211.5004576666667
$ws.Range("N7").Value = 634.5013730000001
$ws.Range("O7").Value = 0.3847969970689267
$ws.Range("P7").Value = 0.3847969970689267
$ws.Range("Q7").Value = 182.4290852590103
$ws.Range("R7").Value = 1641.861767331093
$ws.Range("S7").Value = 0.1949351159494173
$ws.Range("T7").Value = 0.1949351159494173
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8625470000000001
$ws.Range("H8").Value = 2.587641
$ws.Range("I8").Value = 0.5065920925430184
$ws.Range("J8").Value = 0.5065920925430184
$ws.Range("M8").Value = 149.6042426666667
$ws.Range("N8").Value = 448.812728
$ws.Range("O8").Value = 0.2721850532240109
$ws.Range("P8").Value = 0.2721850532240109
$ws.Range("Q8").Value = 129.0406906994053
$ws.Range("R8").Value = 1161.366216294648
$ws.Range("S8").Value = 0.1378867956716845
$ws.Range("T8").Value = 0.1378867956716845
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8625470000000001
$ws.Range("H9").Value = 2.587641
$ws.Range("I9").Value = 0.5065920925430184
$ws.Range("J9").Value = 0.5065920925430184
$ws.Range("M9").Value = 180.1047823333333
$ws.Range("N9").Value = 540.314347
$ws.Range("O9").Value = 0.3276767349073302
$ws.Range("P9").Value = 0.3276767349073302
$ws.Range("Q9").Value = 155.3488396872697
$ws.Range("R9").Value = 1398.139557185427
$ws.Range("S9").Value = 0.1659984428143683
$ws.Range("T9").Value = 0.1659984428143683
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.1658013333333333
$ws.Range("H10").Value = 0.497404
$ws.Range("I10").Value = 0.09737862910630474
$ws.Range("J10").Value = 0.09737862910630474
$ws.Range("M10").Value = 8.432170666666666
$ws.Range("N10").Value = 25.296512
$ws.Range("O10").Value = 0.0153412147997323
$ws.Range("P10").Value = 0.01534121479973231
$ws.Range("Q10").Value = 1.398065139427555
$ws.Range("R10").Value = 12.582586254848
$ws.Range("S10").Value = 0.001493906466023285
$ws.Range("T10").Value = 0.001493906466023285
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.1658013333333333
$ws.Range("H11").Value = 0.497404
$ws.Range("I11").Value = 0.09737862910630474
$ws.Range("J11").Value = 0.09737862910630474
$ws.Range("M11").Value = 211.5004576666667
$ws.Range("N11").Value = 634.5013730000001
$ws.Range("O11").Value = 0.3847969970689267
$ws.Range("P11").Value = 0.3847969970689267
$ws.Range("Q11").Value = 35.06705788174356
$ws.Range("R11").Value = 315.603520935692
$ws.Range("S11").Value = 0.03747100405879484
$ws.Range("T11").Value = 0.03747100405879484
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.1658013333333333
$ws.Range("H12").Value = 0.497404
$ws.Range("I12").Value = 0.09737862910630474
$ws.Range("J12").Value = 0.09737862910630474
$ws.Range("M12").Value = 149.6042426666667
$ws.Range("N12").Value = 448.812728
$ws.Range("O12").Value = 0.2721850532240109
$ws.Range("P12").Value = 0.2721850532240109
$ws.Range("Q12").Value = 24.80458290645689
$ws.Range("R12").Value = 223.241246158112
$ws.Range("S12").Value = 0.02650500734618077
$ws.Range("T12").Value = 0.02650500734618077
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.1658013333333333
$ws.Range("H13").Value = 0.497404
$ws.Range("I13").Value = 0.09737862910630474
$ws.Range("J13").Value = 0.09737862910630474
$ws.Range("M13").Value = 180.1047823333333
$ws.Range("N13").Value = 540.314347
$ws.Range("O13").Value = 0.3276767349073302
$ws.Range("P13").Value = 0.3276767349073302
$ws.Range("Q13").Value = 29.86161305057644
$ws.Range("R13").Value = 268.754517455188
$ws.Range("S13").Value = 0.03190871123530584
$ws.Range("T13").Value = 0.03190871123530584
